$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 60: G1 / Test1 entry for 2025-08-21
$ws.Range("A60").Value = "G1"
$ws.Range("B60").Value = "Test1"
$ws.Range("C60").Value = 45890
$ws.Range("C60").NumberFormat = $ws.Range("C59").NumberFormat
$ws.Range("D60").Value = 0.7568355684346453
$ws.Range("E60").Value = 0
$ws.Range("F60").Value = -0.01

# Row 61: G2 / sedrftgyhuioygtfrd entry for 2025-08-21
$ws.Range("A61").Value = "G2"
$ws.Range("B61").Value = "sedrftgyhuioygtfrd"
$ws.Range("C61").Value = 45890
$ws.Range("C61").NumberFormat = $ws.Range("C59").NumberFormat
$ws.Range("D61").Value = 0.7568355684346453
$ws.Range("E61").Value = 0
$ws.Range("F61").Value = -0.01
